# Fall_2022_Edit_1.05_Students.xlsx -- "2nd Algorithm Part Done" commit
#
# - Adds two new students (rows 101 / 102) to Student_Info with full survey
#   answers, including a Partner_EID column referencing existing students.
# - Flips AM5:AO5 on Student_Info from 0/0/0 to 1/1/1.
# - Moves the active tab / selection from Project_Preferences back to
#   Student_Info (scrolled down near the newly-added rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student_Info")

# ---------------------------------------------------------------------
# AM5:AO5 0 -> 1
# ---------------------------------------------------------------------
$ws.Range("AM5").Value = 1
$ws.Range("AN5").Value = 1
$ws.Range("AO5").Value = 1

# ---------------------------------------------------------------------
# New row 101 (Name100 / EID100)
# ---------------------------------------------------------------------
$ws.Range("A101").Value = "Name100"
$ws.Range("B101").Value = "EID100"
$ws.Range("C101").Value = 3.3
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = "EID037"
$ws.Range("J101").Value = 2
$ws.Range("K101").Value = 4
$ws.Range("L101").Value = 2
$ws.Range("M101").Value = 1
$ws.Range("N101").Value = 1
$ws.Range("O101").Value = 1
$ws.Range("P101").Value = 4
$ws.Range("Q101").Value = 4
$ws.Range("R101").Value = 1
$ws.Range("S101").Value = 4
$ws.Range("T101").Value = 5
$ws.Range("U101").Value = 3
$ws.Range("V101").Value = 4
$ws.Range("W101").Value = 5
$ws.Range("X101").Value = 5
$ws.Range("Y101").Value = 4
$ws.Range("Z101").Value = 5
$ws.Range("AA101").Value = 5
$ws.Range("AB101").Value = 2
$ws.Range("AC101").Value = 3
$ws.Range("AD101").Value = 1
$ws.Range("AE101").Value = 2
$ws.Range("AF101").Value = 3
$ws.Range("AG101").Value = 1
$ws.Range("AH101").Value = 4
$ws.Range("AI101").Value = 2
$ws.Range("AJ101").Value = 1
$ws.Range("AK101").Value = 2
$ws.Range("AL101").Value = 3
$ws.Range("AM101").Value = 2
$ws.Range("AN101").Value = 1
$ws.Range("AO101").Value = 3
$ws.Range("AP101").Value = 1
$ws.Range("AQ101").Value = 1

# ---------------------------------------------------------------------
# New row 102 (Name101 / EID101)
# ---------------------------------------------------------------------
$ws.Range("A102").Value = "Name101"
$ws.Range("B102").Value = "EID101"
$ws.Range("C102").Value = 3.7
$ws.Range("D102").Value = 1
$ws.Range("E102").Value = 1
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = "EID096"
$ws.Range("J102").Value = 1
$ws.Range("K102").Value = 1
$ws.Range("L102").Value = 1
$ws.Range("M102").Value = 1
$ws.Range("N102").Value = 1
$ws.Range("O102").Value = 1
$ws.Range("P102").Value = 1
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = 5
$ws.Range("S102").Value = 1
$ws.Range("T102").Value = 1
$ws.Range("U102").Value = 1
$ws.Range("V102").Value = 5
$ws.Range("W102").Value = 4
$ws.Range("X102").Value = 2
$ws.Range("Y102").Value = 3
$ws.Range("Z102").Value = 3
$ws.Range("AA102").Value = 4
$ws.Range("AB102").Value = 2
$ws.Range("AC102").Value = 4
$ws.Range("AD102").Value = 3
$ws.Range("AE102").Value = 1
$ws.Range("AF102").Value = 5
$ws.Range("AG102").Value = 3
$ws.Range("AH102").Value = 2
$ws.Range("AI102").Value = 1
$ws.Range("AJ102").Value = 4
$ws.Range("AK102").Value = 5
$ws.Range("AL102").Value = 2
$ws.Range("AM102").Value = 3
$ws.Range("AN102").Value = 4
$ws.Range("AO102").Value = 2
$ws.Range("AP102").Value = 1
$ws.Range("AQ102").Value = 3

# ---------------------------------------------------------------------
# View state: bring Student_Info to the front (was Project_Preferences),
# scroll near the bottom of the new data and select AS96, matching the
# author's final cursor position. Project_Preferences loses its
# tabSelected / topLeftCell scroll state as a side effect of switching
# the active sheet.
# ---------------------------------------------------------------------
$ws.Select()
$excel.ActiveWindow.ScrollRow = 93
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("AS96").Select()
